# Apply cryptocurrency price/volume refresh for Thu Feb  1 07:27:15 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell whose style (index 0 / default "Normal") we reuse to restore
# formatting after forcing a numeric-looking value to be stored as text
# (mirrors how the source data keeps these as plain strings, not numbers).
$defaultStyle = $ws.Range("A2").Style

$ws.Range("D2").Value = "42.204.23"
$ws.Range("E2").Value = "  -1.83%  "

$ws.Range("D3").Value = "2.274.75"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.19"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -2.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.11"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -6.00%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -3.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -3.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.25"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -4.61%  "

$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("E12").Value = "  -8.26%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.63"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.74"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "2.627.65"
$ws.Range("E16").Value = "  -2.67%  "

$ws.Range("D17").Value = "2.293.59"
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.778"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  -5.87%  "

$ws.Range("D19").Value = "42.182.44"

$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.45"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -2.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -3.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.65"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -1.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.99"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.98"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -2.97%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  -4.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.86"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  -6.28%  "

$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.81"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.81"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("E32").Value = "  -3.90%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("E34").Value = "  -4.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -2.26%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0690"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -5.32%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -5.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.08"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -8.55%  "

$ws.Range("E39").Value = "  -4.80%  "

$ws.Range("E40").Value = "  -3.31%  "

$ws.Range("E41").Value = "  -3.63%  "

$ws.Range("E42").Value = "  -7.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -6.07%  "

$ws.Range("D44").Value = "1.958.94"
$ws.Range("E44").Value = "  -3.28%  "

$ws.Range("E45").Value = "  -2.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.42"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -7.44%  "

$ws.Range("E47").Value = "  -5.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -4.71%  "

$ws.Range("D49").Value = "2.498.92"
$ws.Range("E49").Value = "  -2.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.03"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -8.02%  "

$ws.Range("E51").Value = "  -5.06%  "
